# Generate Report for Handoff
# Adds a new localization-status row (f3530d67-1e93-4e30-8951-bf0f5929d54c.md)
# to the Overview / zh-cn / de-de sheets, resizes the three tables, and wires
# up the new hyperlinks (with their external-URL relationships).

$wb = $excel.ActiveWorkbook

$fileId = "f3530d67-1e93-4e30-8951-bf0f5929d54c"
$mdName = "$fileId.md"
$mdDisplayOverview = "e2e\$fileId.md"
$zhXlf = "$fileId.c73efe8aeb63621a6146fac2f835ada6f9a47b39.zh-cn.xlf"
$deXlf = "$fileId.c73efe8aeb63621a6146fac2f835ada6f9a47b39.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview" -- new row 10 (columns A:G)
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Range("C10").Value = ".md"
$wsOv.Range("D10").Value = ""
$wsOv.Range("E10").Value = "Ready for handoff"
$wsOv.Range("F10").Value = "Ready for handoff"
$wsOv.Range("G10").Value = "2016-12-15 03:56:48"
$wsOv.Range("G10").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$null = $wsOv.Hyperlinks.Add($wsOv.Range("A10"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3c8e6a129e946c0c7b12f08d6b43e3a63c2af31/e2e/$mdName", $null, $null, $mdName)
$null = $wsOv.Hyperlinks.Add($wsOv.Range("B10"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3c8e6a129e946c0c7b12f08d6b43e3a63c2af31/e2e/$mdName", $null, $null, $mdDisplayOverview)

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G10"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -- new row 10 (columns A:R)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B10").Value = ".md"
$wsZh.Range("C10").Value = "Ready for handoff"
$wsZh.Range("D10").Value = "e2e"
$wsZh.Range("E10").Value = "ht"
$wsZh.Range("F10").Value = "'False"
$wsZh.Range("G10").Value = $zhXlf
$wsZh.Range("H10").Value = "2016-12-15 03:56:34"
$wsZh.Range("H10").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I10").Value = ""
$wsZh.Range("J10").Value = ""
$wsZh.Range("K10").Value = ""
$wsZh.Range("L10").Value = "0001-01-01 00:00:00"
$wsZh.Range("L10").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M10").Value = ""
$wsZh.Range("N10").Value = ""
$wsZh.Range("O10").Value = "'True"
$wsZh.Range("P10").Value = ""
$wsZh.Range("Q10").Value = "'False"
$wsZh.Range("R10").Value = ""

$null = $wsZh.Hyperlinks.Add($wsZh.Range("A10"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c6a2f4f9a6bb1c3eb1ea3e0a27b8f6e1d7b98a2e/e2e/$mdName", $null, $null, $mdName)

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:R10"))

# ---------------------------------------------------------------------------
# Sheet "de-de" -- new row 10 (columns A:R)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B10").Value = ".md"
$wsDe.Range("C10").Value = "Ready for handoff"
$wsDe.Range("D10").Value = "e2e"
$wsDe.Range("E10").Value = "ht"
$wsDe.Range("F10").Value = "'False"
$wsDe.Range("G10").Value = $deXlf
$wsDe.Range("H10").Value = "2016-12-15 03:56:48"
$wsDe.Range("H10").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I10").Value = ""
$wsDe.Range("J10").Value = ""
$wsDe.Range("K10").Value = ""
$wsDe.Range("L10").Value = "0001-01-01 00:00:00"
$wsDe.Range("L10").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M10").Value = ""
$wsDe.Range("N10").Value = ""
$wsDe.Range("O10").Value = "'True"
$wsDe.Range("P10").Value = ""
$wsDe.Range("Q10").Value = "'False"
$wsDe.Range("R10").Value = ""

$null = $wsDe.Hyperlinks.Add($wsDe.Range("A10"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f1b6a7d3c4e8b2a9f5c3d6e7a8b9c0d1e2f3a4b5/e2e/$mdName", $null, $null, $mdName)

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:R10"))
